{"js": "// Recreated Mantel correlograms with Euclidean distances: update the\n// \"Mantel r\" and \"p\" columns of the one_indiv_10km correlogram table.\n//\n// Table layout (1 header row + 7 data rows):\n//   Distance Class (m) | N | Mantel r | p\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst grid = table.values;\nconst header = grid[0].map((h) => String(h).trim().toLowerCase());\nlet rCol = header.indexOf(\"mantel r\");\nlet pCol = header.indexOf(\"p\");\nif (rCol === -1) rCol = 2;\nif (pCol === -1) pCol = 3;\n\n// Old -> new value per data row, in table order (row 1..7, i.e. grid[1..7]).\nconst newValues = [\n  { r: \"0.002\", p: \"0.475\" },\n  { r: \"0.003\", p: \"0.845\" },\n  { r: \"-0.007\", p: \"1\" },\n  { r: \"-0.014\", p: \"0.811\" },\n  { r: \"-0.019\", p: \"0.534\" },\n  { r: \"-0.008\", p: \"1\" },\n  { r: \"0.001\", p: \"1\" },\n];\n\nfor (let i = 0; i < newValues.length; i++) {\n  const row = i + 1; // skip header row\n  if (row >= grid.length) break;\n  table.getCell(row, rCol).value = newValues[i].r;\n  table.getCell(row, pCol).value = newValues[i].p;\n}\n\nawait context.sync();\n", "ps1": "# Recreated Mantel correlograms with Euclidean distances: update the\n# \"Mantel r\" and \"p\" columns of the one_indiv_10km correlogram table.\n#\n# Table layout (1 header row + 7 data rows):\n#   Distance Class (m) | N | Mantel r | p\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# Locate the \"Mantel r\" / \"p\" columns from the header row, falling back to\n# the known positions (3, 4) if the header text doesn't match exactly.\n$rCol = 3\n$pCol = 4\n$colCount = $tbl.Columns.Count\nfor ($c = 1; $c -le $colCount; $c++) {\n    $hdr = $tbl.Cell(1, $c).Range.Text.Trim().TrimEnd([char]7).TrimEnd([char]13).Trim()\n    if ($hdr -eq \"Mantel r\") { $rCol = $c }\n    if ($hdr -eq \"p\") { $pCol = $c }\n}\n\n# Old -> new value per data row, in table order (rows 2..8, header is row 1).\n$newValues = @(\n    @{ r = \"0.002\";  p = \"0.475\" },\n    @{ r = \"0.003\";  p = \"0.845\" },\n    @{ r = \"-0.007\"; p = \"1\" },\n    @{ r = \"-0.014\"; p = \"0.811\" },\n    @{ r = \"-0.019\"; p = \"0.534\" },\n    @{ r = \"-0.008\"; p = \"1\" },\n    @{ r = \"0.001\";  p = \"1\" }\n)\n\nfor ($i = 0; $i -lt $newValues.Count; $i++) {\n    $row = $i + 2  # skip the header row\n    $tbl.Cell($row, $rCol).Range.Text = $newValues[$i].r\n    $tbl.Cell($row, $pCol).Range.Text = $newValues[$i].p\n}\n"}
